# Auto-generated edit script applying the cryptos-list refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.518.18"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3
$ws.Range("D3").Value = "3.116.81"
$ws.Range("E3").Value = "  +0.34%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'526.93"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6
$ws.Range("D6").Value = "'137.29"
$ws.Range("E6").Value = "  -2.79%  "

# Row 8
$ws.Range("D8").Value = "3.113.65"
$ws.Range("E8").Value = "  +0.24%  "

# Row 9
$ws.Range("D9").Value = "'0.447"
$ws.Range("E9").Value = "  +2.72%  "

# Row 10
$ws.Range("D10").Value = "'7.21"
$ws.Range("E10").Value = "  -0.42%  "

# Row 11
$ws.Range("E11").Value = "  -0.79%  "

# Row 12
$ws.Range("D12").Value = "'0.397"
$ws.Range("E12").Value = "  +3.20%  "

# Row 13
$ws.Range("D13").Value = "3.654.16"
$ws.Range("E13").Value = "  +0.38%  "

# Row 14
$ws.Range("E14").Value = "  +2.76%  "

# Row 15
$ws.Range("D15").Value = "'25.38"
$ws.Range("E15").Value = "  -2.93%  "

# Row 16
$ws.Range("E16").Value = "  +0.11%  "

# Row 17
$ws.Range("D17").Value = "57.608.49"
$ws.Range("E17").Value = "  +0.06%  "

# Row 18
$ws.Range("D18").Value = "3.111.43"
$ws.Range("E18").Value = "  +0.33%  "

# Row 19
$ws.Range("D19").Value = "'5.96"
$ws.Range("E19").Value = "  -2.46%  "

# Row 20
$ws.Range("D20").Value = "'12.61"
$ws.Range("E20").Value = "  -1.34%  "

# Row 21
$ws.Range("D21").Value = "'7.92"
$ws.Range("E21").Value = "  -1.81%  "

# Row 22
$ws.Range("D22").Value = "'348.89"
$ws.Range("E22").Value = "  +3.73%  "

# Row 23
$ws.Range("E23").Value = "  -0.64%  "

# Row 24
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").Value = "'68.29"
$ws.Range("E25").Value = "  +2.64%  "

# Row 26
$ws.Range("D26").Value = "'0.505"
$ws.Range("E26").Value = "  -1.42%  "

# Row 27
$ws.Range("E27").Value = "  -0.69%  "

# Row 28
$ws.Range("D28").Value = "'0.994"
$ws.Range("E28").Value = "  -0.79%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0912"
$ws.Range("E29").Value = "  -1.11%  "

# Row 30
$ws.Range("D30").Value = "'7.45"
$ws.Range("E30").Value = "  +3.41%  "

# Row 31
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
$ws.Range("D32").Value = "'1.88"
$ws.Range("E32").Value = "  +0.87%  "

# Row 33
$ws.Range("D33").Value = "'6.05"
$ws.Range("E33").Value = "  -7.28%  "

# Row 34
$ws.Range("D34").Value = "'21.08"
$ws.Range("E34").Value = "  +0.81%  "

# Row 36
$ws.Range("D36").Value = "'4.98"
$ws.Range("E36").Value = "  +7.15%  "

# Row 37
$ws.Range("D37").Value = "'157.90"
$ws.Range("E37").Value = "  +0.60%  "

# Row 38
$ws.Range("D38").Value = "'6.16"
$ws.Range("E38").Value = "  +0.81%  "

# Row 39
$ws.Range("D39").Value = "'26.19"
$ws.Range("E39").Value = "  -3.11%  "

# Row 40
$ws.Range("E40").Value = "  -2.78%  "

# Row 41
$ws.Range("E41").Value = "  +0.83%  "

# Row 42
$ws.Range("D42").Value = "'4.19"
$ws.Range("E42").Value = "  +6.38%  "

# Row 43
$ws.Range("E43").Value = "  +7.22%  "

# Row 44
$ws.Range("D44").Value = "'0.700"
$ws.Range("E44").Value = "  +2.06%  "

# Row 45
$ws.Range("D45").Value = "3.154.17"
$ws.Range("E45").Value = "  +0.24%  "

# Row 46
$ws.Range("D46").Value = "'36.54"
$ws.Range("E46").Value = "  -0.80%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.345.76"
$ws.Range("E47").Value = "  +2.04%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0269"
$ws.Range("E48").Value = "  +3.65%  "

# Row 49
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.06%  "

# Row 50
$ws.Range("D50").Value = "'0.958"
$ws.Range("E50").Value = "  -2.04%  "

# Row 51
$ws.Range("D51").Value = "'6.04"
$ws.Range("E51").Value = "  +0.37%  "
